$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 19 by copying the existing row 19 (OPQA-500_1), which
# already carries the exact cell-style template (s=8/11/20/7/8/5 with the extra
# empty-but-styled G/H/I cells) that the new rows need. Copy+Insert duplicates
# both formatting and values, and shifts the old rows 19-37 down to 21-39.
$ws.Rows.Item(19).Copy()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(19).Copy()
$ws.Rows.Item(19).Insert()

# Row heights for the two new rows
$ws.Rows.Item(19).RowHeight = 32.25
$ws.Rows.Item(20).RowHeight = 32.25

# --- Row 19: OPQA-3176 - Verify that user media image is deleted using Vanilla Profile API ---
$ws.Range("A19").Value2 = "OPQA-3176"
$ws.Range("B19").Value2 = "Verify that user media image is deleted using Vanilla Profile API"
$ws.Range("C19").Value2 = "1PPROFILE"
$ws.Range("D19").Value2 = "/users/user/(SYS_USER2)/media/image"
$ws.Range("E19").Value2 = "DELETE"
$ws.Range("F19").ClearContents()
$ws.Range("I19").Value2 = "OPQA-500"
$ws.Range("J19").Value2 = "status=200"

# --- Row 20: OPQA-3177 - Verify that imageUrl is not returned in Get User Profile API once media image is deleted using Vanilla API ---
$ws.Range("A20").Value2 = "OPQA-3177"
$ws.Range("B20").Value2 = "Verify that imageUrl is not returned in Get User Profile API once media image is deleted using Vanilla API"
$ws.Range("C20").Value2 = "1PPROFILE"
$ws.Range("D20").Value2 = "/users/user/(SYS_USER2)"
$ws.Range("E20").Value2 = "GET"
$ws.Range("F20").ClearContents()
$ws.Range("I20").Value2 = "OPQA-3176"
$ws.Range("J20").Value2 = "status=200||truid=(SYS_USER2)||mediaCategory=image-full"

# Update the top-left cell / selection so the view matches the authored state
# (scrolled so the new rows toward the bottom of the used range are visible).
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("L2:L39").Select()
